# Add a new "schedule" row (row 39) to the active sheet, continuing the
# weekly date / period pattern that already fills rows 1-38.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: next weekly date following 2025/12/12 (row 38)
$ws.Cells.Item(39, 1).Value = "2025/12/19"
# Column B: the date this entry was recorded
$ws.Cells.Item(39, 2).Value = "2026/2/13"
# Column C: description text
$ws.Cells.Item(39, 3).Value = "第89期 第四代寵物"

# Move/select the cell a couple of rows below the newly entered data,
# matching where the cursor ended up after the edit.
[void]$ws.Range("C41").Select()
